$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rating_table")

$ws.Range("B13").Value = 0.07727187447510635
$ws.Range("C13").Value = 0.1925982192678802
$ws.Range("D13").Value = 0.04879821520584188
$ws.Range("E13").Value = 0.3626604749641593
$ws.Range("F13").Value = 0.2353880186521438
$ws.Range("G13").Value = 0.05080275595694478
$ws.Range("H13").Value = 0.03248044147792362
